$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '327.18'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-1.20%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '43.72'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '5.55%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.557'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-2.52%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08085'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-4.04%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '8.658'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-1.80%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '4.274'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-4.98%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.886'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-4.67%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.789'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-5.54%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9356'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.83%'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-5.76%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1893'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-4.24%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09557'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.31%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.04150'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '4.67%'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.40%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001265'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.97%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.005898'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-3.54%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.570'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '3.90%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3487'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.70%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.571'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-6.40%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1361'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.24%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '3.07%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04319'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-2.22%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001237'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.75%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004387'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.37%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001231'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '3.32%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003983'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-0.33%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02657'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-5.98%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05459'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-1.10%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01142'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '27.50%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007683'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-2.74%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1392'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-3.27%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002108'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '1.22%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009236'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006988'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-4.57%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000749'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.33%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003556'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '10.73%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002267'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.66%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002097'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.33%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001997'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.33%'
